# Add season-record columns (Wins, Losses, Ties) to the roster/stats sheet.
# Mirrors the existing "Unnamed: 28" header style for the three new header
# cells, then fills the same record (85-77-0) down every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the style of the last existing header cell
# (AC1) onto the three new header cells so they pick up the same bold /
# bordered / centered formatting, then overwrite their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-58): every row gets the same 2013 Yankees season record.
$lastRow = 58

$ws.Range("AD2:AD$lastRow").Value = 85
$ws.Range("AE2:AE$lastRow").Value = 77
$ws.Range("AF2:AF$lastRow").Value = 0
